# Update cryptos list: refresh Price (D) and Volume(1h) (E) values,
# and fix the ordering of BitcoinSV / MultiversX rows (50 and 51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "42.864.56"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "2.367.04"
$ws.Range("E3").Value = "  +2.14%  "
$ws.Range("E4").Value = "  +0.02%  "
Set-CellText "D5" "301.26"
$ws.Range("E5").Value = "  -0.38%  "
Set-CellText "D6" "95.89"
$ws.Range("E6").Value = "  -0.07%  "
Set-CellText "D7" "0.504"
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("E8").Value = "  -0.08%  "
Set-CellText "D9" "0.495"
$ws.Range("E9").Value = "  -0.05%  "
Set-CellText "D10" "34.09"
$ws.Range("E10").Value = "  -1.18%  "
Set-CellText "D11" "0.0787"
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("E12").Value = "  +2.50%  "
Set-CellText "D13" "18.27"
$ws.Range("E13").Value = "  -3.33%  "
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("D15").Value = "2.742.03"
$ws.Range("E15").Value = "  +2.24%  "
$ws.Range("D16").Value = "2.370.17"
$ws.Range("E16").Value = "  +2.13%  "
Set-CellText "D17" "0.800"
$ws.Range("E17").Value = "  +1.04%  "
$ws.Range("D18").Value = "42.846.79"
$ws.Range("E18").Value = "  +0.11%  "
Set-CellText "D19" "12.13"
$ws.Range("E19").Value = "  -0.21%  "
Set-CellText "D20" "6.30"
$ws.Range("E20").Value = "  +2.23%  "
$ws.Range("D21").Value = "0.0₃0887"
$ws.Range("E21").Value = "  -0.63%  "
Set-CellText "D22" "68.00"
$ws.Range("E22").Value = "  +0.11%  "
Set-CellText "D23" "234.89"
$ws.Range("E23").Value = "  -0.47%  "
Set-CellText "D24" "2.22"
$ws.Range("E24").Value = "  -1.35%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("E26").Value = "  +0.64%  "
Set-CellText "D27" "24.84"
$ws.Range("E27").Value = "  +1.98%  "
Set-CellText "D28" "2.35"
$ws.Range("E28").Value = "  -0.20%  "
Set-CellText "D29" "9.21"
$ws.Range("E29").Value = "  +0.85%  "
Set-CellText "D30" "31.63"
$ws.Range("E30").Value = "  -1.99%  "
$ws.Range("E31").Value = "  -0.01%  "
Set-CellText "D32" "5.04"
$ws.Range("E32").Value = "  +0.81%  "
$ws.Range("E33").Value = "  +4.84%  "
Set-CellText "D34" "17.37"
$ws.Range("E34").Value = "  -2.85%  "
$ws.Range("E35").Value = "  +5.36%  "
$ws.Range("E36").Value = "  +4.35%  "
$ws.Range("E37").Value = "  -2.51%  "
$ws.Range("E38").Value = "  -0.80%  "
Set-CellText "D39" "2.79"
$ws.Range("E39").Value = "  +1.76%  "
Set-CellText "D40" "22.41"
$ws.Range("E40").Value = "  +7.38%  "
$ws.Range("E41").Value = "  -0.62%  "
Set-CellText "D42" "117.77"
$ws.Range("E42").Value = "  -29.15%  "
$ws.Range("D43").Value = "1.937.59"
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("E44").Value = "  +0.47%  "
$ws.Range("E45").Value = "  +1.88%  "
$ws.Range("E46").Value = "  -0.95%  "
Set-CellText "D47" "9.18"
$ws.Range("E47").Value = "  -9.69%  "
$ws.Range("D48").Value = "2.599.15"
$ws.Range("E48").Value = "  +1.95%  "
$ws.Range("E49").Value = "  +2.01%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-CellText "D50" "51.99"
$ws.Range("E50").Value = "  -2.51%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-CellText "D51" "71.97"
$ws.Range("E51").Value = "  -0.16%  "
